$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: the stale yopmail test account is replaced with a live
#     yahoo.com login. It no longer needs its own mailto hyperlink (the
#     same address is already linked from A4), but the cell keeps its
#     existing "Hyperlink" look.
#     NOTE: this engine's Range.Hyperlinks.Delete() always clears every
#     hyperlink on the sheet (there's no working single-link delete), so we
#     clear them all up front and re-create the ones that must survive.
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "testingdev1@yahoo.com"
$ws.Range("B2").Value = "password"

$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:testingdev1@ahoo.com") | Out-Null
$ws.Range("A3").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:testingdev1@yahoo.com") | Out-Null
$ws.Range("A4").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:testulala@yahoo.com") | Out-Null
$ws.Range("A6").Style = "Hyperlink"

# --- New rows for the two additional Sprint 9 test accounts ---
$ws.Range("A7").Value = "testingdev3@yahoo.com"
$ws.Range("B7").Value = "password"
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:testingdev3@yahoo.com") | Out-Null
$ws.Range("A7").Style = "Hyperlink"

$ws.Range("A8").Value = "testingdev22@yahoo.com"
$ws.Hyperlinks.Add($ws.Range("A8"), "mailto:testingdev22@yahoo.com") | Out-Null
$ws.Range("A8").Style = "Hyperlink"

# --- Column widths refreshed to fit the new (longer) email addresses ---
$ws.Columns.Item(1).ColumnWidth = 26.5
$ws.Columns.Item(2).ColumnWidth = 9.6

$ws.Range("F13").Select()
